$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the active selection to B2
$ws.Range("B2").Select()

# Swap the values in columns B and C for rows 1 (header) through 151 (data)
for ($r = 1; $r -le 151; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)
    $tmp = $bCell.Value()
    $bCell.Value = $cCell.Value()
    $cCell.Value = $tmp
}
